$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: labels over the new ratio columns (reuse existing "Flaming"/"Smouldering" strings) ---
$ws.Range("T6").Value = "Flaming"
$ws.Range("U6").Value = "Smouldering"

# --- Rows 20-22: PM2.5 / PM10 / NMOG emission ratios to CO ---
$ws.Range("C25").Value = "DC"
$ws.Range("E25").Value = "Frac Consunm"
$ws.Range("C23").Value = "a"
$ws.Range("D23").Value = "b"
$ws.Range("E23").Value = "c"
$ws.Range("F25").Value = "InvLogit(FracConsump)"
$ws.Range("D25").Value = "Duff Load (kg m2)"

$ws.Range("S20").Value = "PM2.5:CO"
$ws.Range("T20").Formula = "=T13/T`$10"
$ws.Range("U20").Formula = "=U13/U`$10"
$ws.Range("T20:U20").Style = "Normal"

$ws.Range("S21").Value = "PM10:CO"
$ws.Range("T21").Formula = "=T11/T`$10"
$ws.Range("U21").Formula = "=U11/U`$10"
$ws.Range("T21:U21").Style = "Normal"

$ws.Range("S22").Value = "NMOG:CO"
$ws.Range("T22").Formula = "=T12/T`$10"
$ws.Range("U22").Formula = "=U12/U`$10"
$ws.Range("T22:U22").Style = "Normal"

# --- Rows 23-29: duff consumption / InvLogit model loop ---
$ws.Range("C24").Value = 0.33
$ws.Range("D24").Value = -0.17
$ws.Range("E24").Value = -4.8

$ws.Range("C26").Value = 200
$ws.Range("D26").Value = 1
$ws.Range("E26").Formula = "=(C26^`$C`$24)*(D26^`$D`$24)+`$E`$24"
$ws.Range("F26").Formula = "=EXP(E26)/(1+EXP(E26))"

$ws.Range("C27").Value = 200
$ws.Range("D27").Value = 3
$ws.Range("E27").Formula = "=(C27^`$C`$24)*(D27^`$D`$24)+`$E`$24"
$ws.Range("F27").Formula = "=EXP(E27)/(1+EXP(E27))"

$ws.Range("C28").Value = 200
$ws.Range("D28").Value = 5
$ws.Range("E28").Formula = "=(C28^`$C`$24)*(D28^`$D`$24)+`$E`$24"
$ws.Range("F28").Formula = "=EXP(E28)/(1+EXP(E28))"

$ws.Range("D29").Value = 7

# --- view state: scroll sheet over and move the active cell to where the loop landed ---
$ws.Range("U22").Select()
